$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Typography")

# Row 4 ("Default" typography): set Wildcard Characters (G4) and Wildcard Ranges (I4)
# so that TouchGFX generates glyphs for the characters/ranges used at runtime.
$ws.Range("G4").Value = '"+-*/. "'
$ws.Range("I4").Value = "a-z,A-Z"

# Rows 5 and 6 ("Large" and "Small" typographies) are removed from the table
# (columns B through H are cleared out, leaving the rows empty in that range).
$ws.Range("B5:H5").ClearContents()
$ws.Range("B6:H6").ClearContents()
